# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on the per-locale
# handback status sheets, for rows 2 and 5 of each sheet.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-18 00:16:35"
$wsZhCn.Range("H2").Value = "2016-03-18 00:16:54"
$wsZhCn.Range("E5").Value = "2016-03-18 00:16:35"
$wsZhCn.Range("H5").Value = "2016-03-18 00:16:54"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-18 00:16:39"
$wsDeDe.Range("H2").Value = "2016-03-18 00:17:00"
$wsDeDe.Range("E5").Value = "2016-03-18 00:16:39"
$wsDeDe.Range("H5").Value = "2016-03-18 00:17:00"
